# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
#
# Updates the "Metadata" sheet (version bump 5.0.0 -> 6.0.0, new date, publisher
# and jurisdiction info replacing the old duplicated "Contact" rows) and the
# "Elements" sheet (refresh the Short/Definition text on the root Extension row
# to reflect the employee-age-in-years specific wording).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

# The sheet used to have two consecutive rows (10 & 11) both reading
# "Contact" / "No display for ContactDetail". Remove the duplicate (row 11);
# everything below shifts up by one row.
$ws1.Rows.Item(11).Delete()

# Bump the StructureDefinition version and publication date.
$ws1.Range("B3").Value = "6.0.0"
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value, and the old "Contact" row has become a
# "Jurisdiction" row.
$ws1.Range("B9").Value = "Alvearie Team"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# On the Elements sheet, the root Extension row's Short/Definition text is
# updated to describe this specific extension instead of the generic text.
$ws2.Range("K2").Value = "Employee Age In Years"
$ws2.Range("L2").Value = "Age in years at the time of the event for the employee"
